$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Ednra"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 150.0354306666667
$ws.Cells.Item(2, 8).Value = 450.106292
$ws.Cells.Item(2, 9).Value = 0.4152507364956075
$ws.Cells.Item(2, 10).Value = 0.4152507364956075
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 6.433013333333332
$ws.Cells.Item(2, 14).Value = 19.29904
$ws.Cells.Item(2, 15).Value = 0.1097146002786867
$ws.Cells.Item(2, 16).Value = 0.1097146002786867
$ws.Cells.Item(2, 17).Value = 965.1799259510753
$ws.Cells.Item(2, 18).Value = 8686.619333559678
$ws.Cells.Item(2, 19).Value = 0.04555906857004582
$ws.Cells.Item(2, 20).Value = 0.04555906857004582
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Ednra"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 150.0354306666667
$ws.Cells.Item(3, 8).Value = 450.106292
$ws.Cells.Item(3, 9).Value = 0.4152507364956075
$ws.Cells.Item(3, 10).Value = 0.4152507364956075
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 32.709374
$ws.Cells.Item(3, 14).Value = 98.12812199999999
$ws.Cells.Item(3, 15).Value = 0.5578561255548565
$ws.Cells.Item(3, 16).Value = 0.5578561255548566
$ws.Cells.Item(3, 17).Value = 4907.565014927069
$ws.Cells.Item(3, 18).Value = 44168.08513434362
$ws.Cells.Item(3, 19).Value = 0.2316501669952403
$ws.Cells.Item(3, 20).Value = 0.2316501669952403
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Ednra"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 150.0354306666667
$ws.Cells.Item(4, 8).Value = 450.106292
$ws.Cells.Item(4, 9).Value = 0.4152507364956075
$ws.Cells.Item(4, 10).Value = 0.4152507364956075
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.07263
$ws.Cells.Item(4, 14).Value = 0.21789
$ws.Cells.Item(4, 15).Value = 0.001238699658362439
$ws.Cells.Item(4, 16).Value = 0.001238699658362439
$ws.Cells.Item(4, 17).Value = 10.89707332932
$ws.Cells.Item(4, 18).Value = 98.07365996388
$ws.Cells.Item(4, 19).Value = 0.0005143709454318601
$ws.Cells.Item(4, 20).Value = 0.0005143709454318601
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Ednra"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 150.0354306666667
$ws.Cells.Item(5, 8).Value = 450.106292
$ws.Cells.Item(5, 9).Value = 0.4152507364956075
$ws.Cells.Item(5, 10).Value = 0.4152507364956075
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 19.41905066666667
$ws.Cells.Item(5, 14).Value = 58.257152
$ws.Cells.Item(5, 15).Value = 0.3311905745080943
$ws.Cells.Item(5, 16).Value = 0.3311905745080943
$ws.Cells.Item(5, 17).Value = 2913.545629911154
$ws.Cells.Item(5, 18).Value = 26221.91066920038
$ws.Cells.Item(5, 19).Value = 0.1375271299848895
$ws.Cells.Item(5, 20).Value = 0.1375271299848895
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Ednra"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 68.382243
$ws.Cells.Item(6, 8).Value = 205.146729
$ws.Cells.Item(6, 9).Value = 0.1892604742946246
$ws.Cells.Item(6, 10).Value = 0.1892604742946246
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 6.433013333333332
$ws.Cells.Item(6, 14).Value = 19.29904
$ws.Cells.Item(6, 15).Value = 0.1097146002786867
$ws.Cells.Item(6, 16).Value = 0.1097146002786867
$ws.Cells.Item(6, 17).Value = 439.90388098224
$ws.Cells.Item(6, 18).Value = 3959.134928840159
$ws.Cells.Item(6, 19).Value = 0.02076463728578939
$ws.Cells.Item(6, 20).Value = 0.02076463728578939
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Ednra"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 68.382243
$ws.Cells.Item(7, 8).Value = 205.146729
$ws.Cells.Item(7, 9).Value = 0.1892604742946246
$ws.Cells.Item(7, 10).Value = 0.1892604742946246
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 32.709374
$ws.Cells.Item(7, 14).Value = 98.12812199999999
$ws.Cells.Item(7, 15).Value = 0.5578561255548565
$ws.Cells.Item(7, 16).Value = 0.5578561255548566
$ws.Cells.Item(7, 17).Value = 2236.740361245882
$ws.Cells.Item(7, 18).Value = 20130.66325121294
$ws.Cells.Item(7, 19).Value = 0.1055801149106738
$ws.Cells.Item(7, 20).Value = 0.1055801149106738
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Ednra"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 68.382243
$ws.Cells.Item(8, 8).Value = 205.146729
$ws.Cells.Item(8, 9).Value = 0.1892604742946246
$ws.Cells.Item(8, 10).Value = 0.1892604742946246
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.07263
$ws.Cells.Item(8, 14).Value = 0.21789
$ws.Cells.Item(8, 15).Value = 0.001238699658362439
$ws.Cells.Item(8, 16).Value = 0.001238699658362439
$ws.Cells.Item(8, 17).Value = 4.96660230909
$ws.Cells.Item(8, 18).Value = 44.69942078181
$ws.Cells.Item(8, 19).Value = 0.0002344368848502647
$ws.Cells.Item(8, 20).Value = 0.0002344368848502646
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Ednra"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 68.382243
$ws.Cells.Item(9, 8).Value = 205.146729
$ws.Cells.Item(9, 9).Value = 0.1892604742946246
$ws.Cells.Item(9, 10).Value = 0.1892604742946246
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 19.41905066666667
$ws.Cells.Item(9, 14).Value = 58.257152
$ws.Cells.Item(9, 15).Value = 0.3311905745080943
$ws.Cells.Item(9, 16).Value = 0.3311905745080943
$ws.Cells.Item(9, 17).Value = 1327.918241517312
$ws.Cells.Item(9, 18).Value = 11951.26417365581
$ws.Cells.Item(9, 19).Value = 0.06268128521331115
$ws.Cells.Item(9, 20).Value = 0.06268128521331114
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "Ednra"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 104.737245
$ws.Cells.Item(10, 8).Value = 314.211735
$ws.Cells.Item(10, 9).Value = 0.2898796499701289
$ws.Cells.Item(10, 10).Value = 0.2898796499701289
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 6.433013333333332
$ws.Cells.Item(10, 14).Value = 19.29904
$ws.Cells.Item(10, 15).Value = 0.1097146002786867
$ws.Cells.Item(10, 16).Value = 0.1097146002786867
$ws.Cells.Item(10, 17).Value = 673.7760935815999
$ws.Cells.Item(10, 18).Value = 6063.984842234398
$ws.Cells.Item(10, 19).Value = 0.0318040299253983
$ws.Cells.Item(10, 20).Value = 0.0318040299253983
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Gnai2"
$ws.Cells.Item(11, 3).Value = "Ednra"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 104.737245
$ws.Cells.Item(11, 8).Value = 314.211735
$ws.Cells.Item(11, 9).Value = 0.2898796499701289
$ws.Cells.Item(11, 10).Value = 0.2898796499701289
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 32.709374
$ws.Cells.Item(11, 14).Value = 98.12812199999999
$ws.Cells.Item(11, 15).Value = 0.5578561255548565
$ws.Cells.Item(11, 16).Value = 0.5578561255548566
$ws.Cells.Item(11, 17).Value = 3425.889718434629
$ws.Cells.Item(11, 18).Value = 30833.00746591167
$ws.Cells.Item(11, 19).Value = 0.1617111384095341
$ws.Cells.Item(11, 20).Value = 0.1617111384095341
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Gnai2"
$ws.Cells.Item(12, 3).Value = "Ednra"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 104.737245
$ws.Cells.Item(12, 8).Value = 314.211735
$ws.Cells.Item(12, 9).Value = 0.2898796499701289
$ws.Cells.Item(12, 10).Value = 0.2898796499701289
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.07263
$ws.Cells.Item(12, 14).Value = 0.21789
$ws.Cells.Item(12, 15).Value = 0.001238699658362439
$ws.Cells.Item(12, 16).Value = 0.001238699658362439
$ws.Cells.Item(12, 17).Value = 7.607066104349999
$ws.Cells.Item(12, 18).Value = 68.46359493915
$ws.Cells.Item(12, 19).Value = 0.000359073823384222
$ws.Cells.Item(12, 20).Value = 0.000359073823384222
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Gnai2"
$ws.Cells.Item(13, 3).Value = "Ednra"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 104.737245
$ws.Cells.Item(13, 8).Value = 314.211735
$ws.Cells.Item(13, 9).Value = 0.2898796499701289
$ws.Cells.Item(13, 10).Value = 0.2898796499701289
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 19.41905066666667
$ws.Cells.Item(13, 14).Value = 58.257152
$ws.Cells.Item(13, 15).Value = 0.3311905745080943
$ws.Cells.Item(13, 16).Value = 0.3311905745080943
$ws.Cells.Item(13, 17).Value = 2033.89786734208
$ws.Cells.Item(13, 18).Value = 18305.08080607872
$ws.Cells.Item(13, 19).Value = 0.09600540781181229
$ws.Cells.Item(13, 20).Value = 0.09600540781181227
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Gnai2"
$ws.Cells.Item(14, 3).Value = "Ednra"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 38.15794
$ws.Cells.Item(14, 8).Value = 114.47382
$ws.Cells.Item(14, 9).Value = 0.105609139239639
$ws.Cells.Item(14, 10).Value = 0.105609139239639
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 6.433013333333332
$ws.Cells.Item(14, 14).Value = 19.29904
$ws.Cells.Item(14, 15).Value = 0.1097146002786867
$ws.Cells.Item(14, 16).Value = 0.1097146002786867
$ws.Cells.Item(14, 17).Value = 245.4705367925333
$ws.Cells.Item(14, 18).Value = 2209.2348311328
$ws.Cells.Item(14, 19).Value = 0.01158686449745316
$ws.Cells.Item(14, 20).Value = 0.01158686449745315
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Gnai2"
$ws.Cells.Item(15, 3).Value = "Ednra"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 38.15794
$ws.Cells.Item(15, 8).Value = 114.47382
$ws.Cells.Item(15, 9).Value = 0.105609139239639
$ws.Cells.Item(15, 10).Value = 0.105609139239639
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 32.709374
$ws.Cells.Item(15, 14).Value = 98.12812199999999
$ws.Cells.Item(15, 15).Value = 0.5578561255548565
$ws.Cells.Item(15, 16).Value = 0.5578561255548566
$ws.Cells.Item(15, 17).Value = 1248.12233052956
$ws.Cells.Item(15, 18).Value = 11233.10097476604
$ws.Cells.Item(15, 19).Value = 0.05891470523940839
$ws.Cells.Item(15, 20).Value = 0.05891470523940838
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Gnai2"
$ws.Cells.Item(16, 3).Value = "Ednra"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 38.15794
$ws.Cells.Item(16, 8).Value = 114.47382
$ws.Cells.Item(16, 9).Value = 0.105609139239639
$ws.Cells.Item(16, 10).Value = 0.105609139239639
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.07263
$ws.Cells.Item(16, 14).Value = 0.21789
$ws.Cells.Item(16, 15).Value = 0.001238699658362439
$ws.Cells.Item(16, 16).Value = 0.001238699658362439
$ws.Cells.Item(16, 17).Value = 2.7714111822
$ws.Cells.Item(16, 18).Value = 24.9427006398
$ws.Cells.Item(16, 19).Value = 0.0001308180046960921
$ws.Cells.Item(16, 20).Value = 0.000130818004696092
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Gnai2"
$ws.Cells.Item(17, 3).Value = "Ednra"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 38.15794
$ws.Cells.Item(17, 8).Value = 114.47382
$ws.Cells.Item(17, 9).Value = 0.105609139239639
$ws.Cells.Item(17, 10).Value = 0.105609139239639
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 19.41905066666667
$ws.Cells.Item(17, 14).Value = 58.257152
$ws.Cells.Item(17, 15).Value = 0.3311905745080943
$ws.Cells.Item(17, 16).Value = 0.3311905745080943
$ws.Cells.Item(17, 17).Value = 740.9909701956268
$ws.Cells.Item(17, 18).Value = 6668.91873176064
$ws.Cells.Item(17, 19).Value = 0.03497675149808137
$ws.Cells.Item(17, 20).Value = 0.03497675149808136
